$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 9 (B9 changes, A9/C9 stay the same)
$ws.Range("B9").Value = 230

# Update row 10 values to merge what used to be rows 10-13
$ws.Range("A10").Value = 230
$ws.Range("B10").Value = 255
$ws.Range("C10").Value = 93

# Update row 11 values to match what used to be row 14
$ws.Range("A11").Value = 255
$ws.Range("B11").Value = 360
$ws.Range("C11").Value = 0

# Delete the now-obsolete rows 12:14 (shifts rows up, clears values/dimension)
$ws.Range("A12:C14").Delete()

# Move the active selection to B16, matching the post-edit workbook state
$ws.Range("B16").Select()
